$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.606.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.15%  "

$ws.Range("D3").Value = "'1.750.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.74%  "

$ws.Range("D4").Value = "'0.9959"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'247.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.84%  "

$ws.Range("D6").Value = "'0.9972"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "'0.4804"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.2719"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.35%  "

$ws.Range("D9").Value = "'0.06252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("D10").Value = "'1.737.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.07%  "

$ws.Range("D11").Value = "'0.07116"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").Value = "'15.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.97%  "

$ws.Range("D13").Value = "'0.6235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.34%  "

$ws.Range("D14").Value = "'4.500"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").Value = "'77.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("D16").Value = "'0.9971"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "'26.602.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.21%  "

$ws.Range("D18").Value = "'0.9962"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").Value = "'0.000006896"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.80%  "

$ws.Range("D20").Value = "'11.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.47%  "

$ws.Range("D21").Value = "'1.958.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.14%  "

$ws.Range("D22").Value = "'4.616"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.60%  "

$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "'5.353"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").Value = "'136.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").Value = "'15.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.62%  "

$ws.Range("D27").Value = "'1.841"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.81%  "

$ws.Range("D28").Value = "'1.411"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "'107.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.78%  "

$ws.Range("D30").Value = "'4.016"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("D31").Value = "'3.773"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("D32").Value = "'0.07900"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").Value = "'0.04579"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.46%  "

$ws.Range("D34").Value = "'2.610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").Value = "'1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.14%  "

$ws.Range("D36").Value = "'0.6330"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("D37").Value = "'0.9541"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.95%  "

$ws.Range("D38").Value = "'114.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.82%  "

$ws.Range("D39").Value = "'2.491"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.98%  "

$ws.Range("D40").Value = "'1.974"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.87%  "

$ws.Range("D41").Value = "'1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("D42").Value = "'5.751"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.84%  "

$ws.Range("D43").Value = "'0.01508"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("D44").Value = "'0.3914"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.787"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.12%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1210"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.14%  "

$ws.Range("D47").Value = "'0.05326"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "

$ws.Range("D48").Value = "'8.033"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.19%  "

$ws.Range("D49").Value = "'30.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.41%  "

$ws.Range("D50").Value = "'0.3452"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.26%  "

$ws.Range("D51").Value = "'1.237"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.15%  "
